$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.074.13"
$ws.Range("E2").Value = "  +7.31%  "
$ws.Range("D3").Value = "2.585.36"
$ws.Range("E3").Value = "  +9.52%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +23.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.997"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "2.585.80"
$ws.Range("E9").Value = "  +9.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.70%  "
$ws.Range("E11").Value = "  +6.14%  "
$ws.Range("E12").Value = "  +6.23%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "3.037.78"
$ws.Range("E14").Value = "  +9.43%  "
$ws.Range("D15").Value = "59.068.24"
$ws.Range("E15").Value = "  +6.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.02%  "
$ws.Range("E17").Value = "  +5.09%  "
$ws.Range("D18").Value = "2.586.51"
$ws.Range("E18").Value = "  +9.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.55%  "
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.86%  "
$ws.Range("E25").Value = "  +5.48%  "
$ws.Range("E26").Value = "  +8.26%  "
$ws.Range("D27").Value = "2.700.94"
$ws.Range("E27").Value = "  +9.27%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0825"
$ws.Range("E29").Value = "  +9.12%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.50%  "
$ws.Range("E37").Value = "  +9.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.80%  "
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "289.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.101"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.625"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.65%  "
$ws.Range("E45").Value = "  +5.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.995"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0236"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.13%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.727"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.38%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.50%  "
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "
